# Split the paragraph that currently holds "New line 3" (and the
# trailing _GoBack bookmark) so that a new paragraph containing
# "New line 4" is inserted right after "New line 3", leaving the
# bookmark alone in its own (now empty) trailing paragraph.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "New line 3",  # Find what
    $false,        # MatchCase
    $false,        # MatchWholeWord
    $false,        # MatchWildcards
    $false,        # MatchSoundsLike
    $false,        # MatchAllWordForms
    $true,         # Forward
    1,             # Wrap (wdFindContinue)
    $false,        # Format
    "New line 3^pNew line 4^p",  # Replace with
    2              # Replace (wdReplaceAll)
)
